$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-45 (columns A-I), reflecting the re-ordered species list,
# updated weights/counts, the two rows reassigned from Gear=2-RAP to Gear=1-RAP,
# the new RF value for the "no-RF" group, and four new trailing rows (42-45).
$data = @(
    @('SOLEMON2024','ITA17','25','1-RAP','Arnoglossus laterna','ARNOLAT',0.012,1,1),
    @('SOLEMON2024','ITA17','25','1-RAP','Bolinus brandaris','MUREBRA',1.232,396,1),
    @('SOLEMON2024','ITA17','25','1-RAP','Callinectes sapidus','CALLSAP',0.191,1,1),
    @('SOLEMON2024','ITA17','25','1-RAP','Chelidonichthys lucernus','TRIGLUC',0.168,1,1),
    @('SOLEMON2024','ITA17','25','1-RAP','Gobius niger','GOBINIG',0.027,1,1),
    @('SOLEMON2024','ITA17','25','1-RAP','Hexaplex trunculus','HEXATRU',1.801,194,1),
    @('SOLEMON2024','ITA17','25','1-RAP','Melicertus kerathurus','MELIKER',0.866,46,1),
    @('SOLEMON2024','ITA17','25','1-RAP','Ostrea edulis','OSTREDU',0.056,2,1),
    @('SOLEMON2024','ITA17','25','1-RAP','Raja asterias','RAJAAST',0.804,1,1),
    @('SOLEMON2024','ITA17','25','1-RAP','Scophthalmus rhombus','SCOHRHO',0.063,1,1),
    @('SOLEMON2024','ITA17','25','1-RAP','Solea solea','SOLEVUL',1.287,8,1),
    @('SOLEMON2024','ITA17','25','1-RAP','Squilla mantis','SQUIMAN',0.201,5,1),
    @('SOLEMON2024','ITA17','25','1-RAP','Tonna galea','TONNGAL',0.222,2,1),
    @('SOLEMON2024','ITA17','25','2-RAP','Bolinus brandaris','MUREBRA',1.12,360,1),
    @('SOLEMON2024','ITA17','25','2-RAP','Chelidonichthys lucernus','TRIGLUC',0.844,7,1),
    @('SOLEMON2024','ITA17','25','2-RAP','Gobius niger','GOBINIG',0.021,2,1),
    @('SOLEMON2024','ITA17','25','2-RAP','Hexaplex trunculus','HEXATRU',1.637,176,1),
    @('SOLEMON2024','ITA17','25','2-RAP','Melicertus kerathurus','MELIKER',1.253,71,1),
    @('SOLEMON2024','ITA17','25','2-RAP','Ostrea edulis','OSTREDU',0.121,1,1),
    @('SOLEMON2024','ITA17','25','2-RAP','Penaeus aztecus','PENAZTC',0.057,1,1),
    @('SOLEMON2024','ITA17','25','2-RAP','Scophthalmus rhombus','SCOHRHO',0.188,1,1),
    @('SOLEMON2024','ITA17','25','2-RAP','Sepia officinalis','SEPIOFF',0.03,1,1),
    @('SOLEMON2024','ITA17','25','2-RAP','Solea solea','SOLEVUL',1.371,10,1),
    @('SOLEMON2024','ITA17','25','2-RAP','Squilla mantis','SQUIMAN',0.114,2,1),
    @('SOLEMON2024','ITA17','25','2-RAP','Actiniaria nd','ACTINND',0.001,1,4.051733333333333),
    @('SOLEMON2024','ITA17','25','2-RAP','Anadara kagoshimensis','ANADKAG',0.053,6,4.051733333333333),
    @('SOLEMON2024','ITA17','25','2-RAP','Aporrhais pespelecani','APORPES',0.668,103,4.051733333333333),
    @('SOLEMON2024','ITA17','25','2-RAP','Ascidiella aspersa','ASCIASP',0.008,1,4.051733333333333),
    @('SOLEMON2024','ITA17','25','2-RAP','Astropecten irregularis','ASTRIRR',0.263,85,4.051733333333333),
    @('SOLEMON2024','ITA17','25','2-RAP','Biological discard','BIOLDIS',0.077,-1,4.051733333333333),
    @('SOLEMON2024','ITA17','25','2-RAP','Bolinus brandaris','MUREBRA',1.691,172,4.051733333333333),
    @('SOLEMON2024','ITA17','25','2-RAP','Eggs of Murex','EGGSMUR',0.014,-1,4.051733333333333),
    @('SOLEMON2024','ITA17','25','2-RAP','Eggs of Raja sp','EGGSRAJ',0.001,1,4.051733333333333),
    @('SOLEMON2024','ITA17','25','2-RAP','Goneplax rhomboides','GONERHO',0.067,11,4.051733333333333),
    @('SOLEMON2024','ITA17','25','2-RAP','Hexaplex trunculus','HEXATRU',0.762,28,4.051733333333333),
    @('SOLEMON2024','ITA17','25','2-RAP','Liocarcinus depurator','LIOCDEP',0.235,24,4.051733333333333),
    @('SOLEMON2024','ITA17','25','2-RAP','Medorippe lanata','MEDOLAN',0.083,10,4.051733333333333),
    @('SOLEMON2024','ITA17','25','2-RAP','Nassarius lima','NASSLIM',0.005,3,4.051733333333333),
    @('SOLEMON2024','ITA17','25','2-RAP','Ophiura ophiura','OPHIOPH',0.018,15,4.051733333333333),
    @('SOLEMON2024','ITA17','25','2-RAP','Schizaster canaliferus','SCHICAN',0.03,1,4.051733333333333),
    @('SOLEMON2024','ITA17','25','2-RAP','Shells NA','SHELLS',0.52,-1,4.051733333333333),
    @('SOLEMON2024','ITA17','25','2-RAP','Trachythyone elongata','TRACELO',0.003,2,4.051733333333333),
    @('SOLEMON2024','ITA17','25','2-RAP','Turritella communis','TURRCOM',0.021,23,4.051733333333333),
    @('SOLEMON2024','ITA17','25','2-RAP','Wood NA','WOOD',0.22,-1,4.051733333333333)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
}
